$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update CHRT_CNT (column J) values per the data correction:
# J2: 1 -> 0
# J8: 0 -> 1
# J14: 0 -> 1
$ws.Range("J2").Value = 0
$ws.Range("J8").Value = 1
$ws.Range("J14").Value = 1

# Force a full recalculation on load, matching calcPr/fullCalcOnLoad="true"
$wb.ForceFullCalculation = $true
